$wb = $excel.ActiveWorkbook

# --- location sheet: tweak two verbosedescription values -------------------
$wsLocation = $wb.Worksheets.Item("location")

# Spooky House: add a trailing period to the verbose description
$wsLocation.Range("H3").Value = "A spooooooooooky house, run by the ever lovely and elderly Dolores, known for her sweet sweet herbs."

# Main Street: rewrite the verbose description
$wsLocation.Range("H6").Value = "Dusty Main Street. A tumbleweed blows slowly by. (this is meant to be verbose)"

# --- command sheet: add a TRAVEL command ------------------------------------
$wsCommand = $wb.Worksheets.Item("command")

# Copy formatting from the existing MOVE row down onto the new row first so
# the new cells pick up the same style as the rest of the table.
$wsCommand.Range("A2:B2").Copy()
$wsCommand.Range("A3:B3").PasteSpecial(-4122)

$wsCommand.Range("A3").Value = "TRAVEL(to: location id)"
$wsCommand.Range("B3").Value = "travel to a place"

# --- item sheet: rename headers + add three new items -----------------------
$wsItem = $wb.Worksheets.Item("item")

# Headers "itemname (string)" / "itemdescription (string)" become the
# generic "name (string)" / "description (string)" (matching other sheets).
$wsItem.Range("B1").Value = "name (string)"
$wsItem.Range("C1").Value = "description (string)"

# Copy formatting from an existing data row onto the three new rows.
$wsItem.Range("A2:E2").Copy()
$wsItem.Range("A4:E6").PasteSpecial(-4122)

$wsItem.Range("A4").Value = 3
$wsItem.Range("B4").Value = "William's Gun"
$wsItem.Range("C4").Value = "Nothing special, it's a gun"
$wsItem.Range("D4").Value = 1
$wsItem.Range("E4").Value = $false

$wsItem.Range("A5").Value = 4
$wsItem.Range("B5").Value = "Key"
$wsItem.Range("C5").Value = "Key to the jail cell"
$wsItem.Range("D5").Value = 1
$wsItem.Range("E5").Value = $true

$wsItem.Range("A6").Value = 5
$wsItem.Range("B6").Value = "Door"
$wsItem.Range("C6").Value = "It's a door"
$wsItem.Range("D6").Value = 1
$wsItem.Range("E6").Value = $true
